$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I3").Value = 11
$ws.Range("J3").Value = 5.5
$ws.Range("P3").Value = 2.62
$ws.Range("U3").Value = 2.06
$ws.Range("AC3").Value = 15
$ws.Range("N4").Value = 5.2
$ws.Range("P4").Value = 2.42
$ws.Range("R4").Value = 1.58
$ws.Range("U4").Value = 2.3
$ws.Range("X4").Value = 23
$ws.Range("AA4").Value = 16.5
$ws.Range("AE4").Value = 16
$ws.Range("AG4").Value = 22
$ws.Range("AI4").Value = 30
$ws.Range("AK4").Value = 70
$ws.Range("AO4").Value = 7.2
$ws.Range("S5").Value = 2.68
$ws.Range("X5").Value = 26
$ws.Range("AL5").Value = 32
$ws.Range("AN5").Value = 18
$ws.Range("G6").Value = 9
$ws.Range("N6").Value = 6.2
$ws.Range("R6").Value = 1.69
$ws.Range("S6").Value = 2.36
$ws.Range("T6").Value = 1.79
$ws.Range("AA6").Value = 13.5
$ws.Range("AC6").Value = 16
$ws.Range("AH6").Value = 25
$ws.Range("AM6").Value = 120
$ws.Range("G7").Value = 2.38
$ws.Range("H7").Value = 3.25
$ws.Range("I7").Value = 3.35
$ws.Range("J7").Value = 3.65
$ws.Range("K7").Value = 3.7
$ws.Range("R7").Value = 1.47
$ws.Range("S7").Value = 2.98
$ws.Range("U7").Value = 2.42
$ws.Range("F8").Value = 1.72
$ws.Range("H8").Value = 5
$ws.Range("I8").Value = 5.3
$ws.Range("N8").Value = 5.8
$ws.Range("U8").Value = 2.44
$ws.Range("AD8").Value = 22
$ws.Range("AI8").Value = 1000
$ws.Range("AJ8").Value = 19
$ws.Range("N9").Value = 5.9
$ws.Range("P9").Value = 2.66
$ws.Range("Q9").Value = 1.56
$ws.Range("S9").Value = 2.4
$ws.Range("T9").Value = 1.8
$ws.Range("G10").Value = 1.36
$ws.Range("H10").Value = 9
$ws.Range("X10").Value = 48
$ws.Range("AA10").Value = 300
$ws.Range("F11").Value = 2.36
$ws.Range("X11").Value = 16
$ws.Range("AG11").Value = 11.5
$ws.Range("AH11").Value = 18
$ws.Range("AJ11").Value = 32
$ws.Range("AL11").Value = 40
$ws.Range("I12").Value = 3.8
$ws.Range("Q12").Value = 1.94
